# BRVM recommandations / Top_YTD automatic refresh (GitHub Actions data pull)
# Overwrites the "Recommandations" and "Top_YTD" sheets of the active workbook
# with the latest pulled figures, and drops the 5 tickers that fell out of the
# Top_YTD/Recommandations ranking window.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Recommandations" sheet
# ---------------------------------------------------------------------------
$wsReco = $wb.Worksheets.Item("Recommandations")

# Latest values for rows 2-42 (Titre, Jours en Hausse, Jours en Baisse,
# Variation Totale (%), Derniere Variation (%), Recommandation, Strategie)
$recoData = @(
    @("BRVM - SERVICES PUBLICS", 0, 8, 3340.07, 105.86, "🟡 Observer", "➖ Neutre"),
    @("NEI-CEDA CI", 0, 4, 2715, 690, "🟡 Observer", "➖ Neutre"),
    @("AIR LIQUIDE CI", 0, 4, 2675, 625, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AUTRES SECTEURS", 0, 4, 2295.06, 592.21, "🟡 Observer", "➖ Neutre"),
    @("BRVM - DISTRIBUTION", 0, 4, 1757.58, 460.26, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AGRICULTURE", 0, 4, 1443.32, 364.03, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TRANSPORT", 0, 4, 1439.27, 356.15, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 4, 567.04, 152.85, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRESTIGE", 0, 4, 547.92, 136.73, "🟡 Observer", "➖ Neutre"),
    @("BRVM - FINANCES", 0, 4, 543.11, 135.67, "🟡 Observer", "➖ Neutre"),
    @("BRVM - SERVICES FINANCIERS", 0, 4, 533.77, 133.34, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIELS", 0, 4, 494.99, 126.36, "🟡 Observer", "➖ Neutre"),
    @("BRVM - ENERGIE", 0, 4, 432.8, 109.2, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE                  (**)", 0, 2, 414.76, 207.68, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL                    (**)", 0, 2, 387.23, 193.83, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TELECOMMUNICATIONS", 0, 4, 380.89, 94.48, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE         (**)", 0, 2, 372.63, 187.07, "🟡 Observer", "➖ Neutre"),
    @("FILTISAC CI (FTSC)", 4, 0, 29.59, 7.4, "🟢 Achat", "✅ Renforcer"),
    @("TRACTAFRIC MOTORS CI (PRSC)", 3, 0, 22.19, 7.39, "🟢 Achat", "✅ Renforcer"),
    @("CFAO MOTORS CI (CFAC)", 3, 1, 16.2, 7.35, "🟢 Achat", "✅ Renforcer"),
    @("BERNABE CI (BNBC)", 2, 0, 14.01, 7.5, "🟡 Observer", "➖ Neutre"),
    @("SICABLE CI (CABC)", 1, 0, 5.2, 5.2, "🟡 Observer", "➖ Neutre"),
    @("VIVO ENERGY CI (SHEC)", 1, 0, 2.92, 2.92, "🟡 Observer", "➖ Neutre"),
    @("ONATEL BF (ONTBF)", 1, 0, 2.61, 2.61, "🟡 Observer", "➖ Neutre"),
    @("NSIA BANQUE COTE D'IVOIRE (NSBC)", 1, 1, 0.49, -2.01, "🟡 Observer", "👀 À surveiller"),
    @("SUCRIVOIRE (SCRC)", 1, 2, 0.43, 6.78, "🟡 Observer", "👀 À surveiller"),
    @("ORAGROUP TOGO (ORGT)", 1, 1, 0.03, -7.38, "🟡 Observer", "👀 À surveiller"),
    @("TOTAL", 0, 4, 0, 0, "🟡 Observer", "➖ Neutre"),
    @("UNIWAX CI (UNXC)", 1, 1, -0.25, 7.23, "🟡 Observer", "👀 À surveiller"),
    @("SAPH CI (SPHC)", 1, 1, -0.32, 5.94, "🟡 Observer", "👀 À surveiller"),
    @("BANK OF AFRICA ML (BOAM)", 0, 1, -1.15, -1.15, "🟡 Observer", "➖ Neutre"),
    @("SOGB CI (SOGC)", 0, 1, -1.18, -1.18, "🟡 Observer", "➖ Neutre"),
    @("AFRICA GLOBAL LOGISTICS CI (SDSC)", 0, 1, -1.68, -1.68, "🟡 Observer", "➖ Neutre"),
    @("TOTALENERGIES MARKETING CI (TTLC)", 0, 1, -2.04, -2.04, "🟡 Observer", "➖ Neutre"),
    @("ECOBANK COTE D''IVOIRE (ECOC)", 0, 1, -2.23, -2.23, "🟡 Observer", "➖ Neutre"),
    @("SOLIBRA CI (SLBC)", 0, 1, -2.75, -2.75, "🟡 Observer", "➖ Neutre"),
    @("NEI-CEDA CI (NEIC)", 0, 1, -2.94, -2.94, "🟡 Observer", "➖ Neutre"),
    @("ORANGE COTE D'IVOIRE (ORAC)", 0, 1, -3.34, -3.34, "🟡 Observer", "➖ Neutre"),
    @("SETAO CI (STAC)", 0, 2, -4.81, -2.91, "🟡 Observer", "➖ Neutre"),
    @("AIR LIQUIDE CI (SIVC)", 0, 1, -7.35, -7.35, "🟡 Observer", "➖ Neutre"),
    @("LOTERIE NATIONALE DU BENIN (LNBB)", 0, 1, -7.45, -7.45, "🟡 Observer", "➖ Neutre")
)

$r = 2
foreach ($row in $recoData) {
    $wsReco.Cells.Item($r, 1).Value = $row[0]
    $wsReco.Cells.Item($r, 2).Value = $row[1]
    $wsReco.Cells.Item($r, 3).Value = $row[2]
    $wsReco.Cells.Item($r, 4).Value = $row[3]
    $wsReco.Cells.Item($r, 5).Value = $row[4]
    $wsReco.Cells.Item($r, 6).Value = $row[5]
    $wsReco.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

# The refreshed dataset is 5 rows shorter (47 -> 42): drop the now-unused
# trailing rows so the sheet's used range shrinks to A1:G42.
$wsReco.Range("A43:G47").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 2) "Top_YTD" sheet
# ---------------------------------------------------------------------------
$wsYtd = $wb.Worksheets.Item("Top_YTD")

# Latest "Progression YTD (%)" figures for rows 2-11 (labels unchanged)
$ytdData = @(
    @("BRVM - SERVICES PUBLICS", 8393274.06),
    @("NEI-CEDA CI", 367524.92),
    @("AIR LIQUIDE CI", 348361.1),
    @("BRVM - AUTRES SECTEURS", 205870.5),
    @("BRVM - DISTRIBUTION", 84428.3),
    @("BRVM - AGRICULTURE", 44991.17),
    @("BRVM - TRANSPORT", 44599.45),
    @("BRVM - CONSOMMATION DISCRETIONNAIRE", 3309.98),
    @("BRVM-PRESTIGE", 3053.88),
    @("BRVM - FINANCES", 2990.36)
)

$r = 2
foreach ($row in $ytdData) {
    $wsYtd.Cells.Item($r, 1).Value = $row[0]
    $wsYtd.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}
